# Word report template: split each "$comments"/"$resolution" placeholder
# paragraph into a pair of paragraphs so the generated report can carry an
# extra image alongside the existing text block. End state (in order):
#   $content  $image_content  $comments  $image_comment  $resolution  $image_resolution

$d = $word.ActiveDocument

function Get-ParaByText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq ($text + "`r")) {
            return $p
        }
    }
    return $null
}

# 1) The paragraph right after "$content" used to hold "$comments" — that
#    slot becomes the new "$image_content" placeholder.
$p = Get-ParaByText("`$comments")
$p.Range.Text = "`$image_content"

# 2) Re-insert "$comments" as its own paragraph right after, followed by a
#    brand new "$image_comment" placeholder paragraph.
$p = Get-ParaByText("`$image_content")
$p.Range.InsertParagraphAfter()
$p.Next().Range.Text = "`$comments"

$p = Get-ParaByText("`$comments")
$p.Range.InsertParagraphAfter()
$p.Next().Range.Text = "`$image_comment"

# 3) After the untouched "$resolution" paragraph, add an "$image_resolution"
#    placeholder paragraph.
$p = Get-ParaByText("`$resolution")
$p.Range.InsertParagraphAfter()
$p.Next().Range.Text = "`$image_resolution"

Write-Output "edit.ps1 applied"
